# se crea punto 42
# Update reporting period (Q1 2023 -> Q2 2023) on "Reporte de Formatos" sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

$ws.Range("B8").Value = 45017
$ws.Range("C8").Value = 45107

$ws.Activate()
$ws.Range("C16").Select()
